# Insert a new weekly "Granada" (pomegranate) price record for
# "Vega Modelo de Temuco" ahead of the existing row 52, shifting every
# subsequent record down by one row (old row 52 -> 53, ..., old row 157 -> 158).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 52..157 down to 53..158, opening up a blank row 52.
$ws.Rows.Item(52).EntireRow.Insert()

# Populate the newly inserted row 52 with the new record.
$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "Vega Modelo de Temuco"
$ws.Range("C52").Value = "La Araucanía"
$ws.Range("D52").Value = 44791
$ws.Range("E52").Value = 9
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100104
$ws.Range("H52").Value = "Frutos de pepita"
$ws.Range("I52").Value = 100104001
$ws.Range("J52").Value = "Granada"
$ws.Range("K52").Value = "Wonderfull"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 400
$ws.Range("N52").Value = 14000
$ws.Range("O52").Value = 15000
$ws.Range("P52").Value = 14500
$ws.Range("Q52").Value = '$/bandeja 10 kilos granel'
$ws.Range("R52").Value = "Provincia de Limarí"
$ws.Range("S52").Value = 1450
$ws.Range("T52").Value = 10
